# Patron.xlsx edit
#
# The data rows for each "Orden" (1-24, rows 6-29) are collapsed/hidden,
# keeping the summary rows 1-5 and the totals rows 25-27 (sheet rows 30-32)
# visible. The user's final selection lands on row 29 (the whole row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide the detail rows 6 through 29 (inclusive); rows 30-32 stay visible.
$ws.Range("A6:A29").EntireRow.Hidden = $true

# Leave the last action as "select row 29" (matches the saved selection
# state A29:XFD29 in the target workbook).
$ws.Range("A29:XFD29").Select()
